$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A7 is an empty-string text cell (matches the other rows' "ملاحظات" column,
# which are all blank). A leading apostrophe forces Excel to store it as
# text (empty string) instead of clearing the cell outright; resetting the
# style back to Normal afterwards avoids leaving a stray number-format.
$ws.Cells.Item(7, 1).Value = "'"
$ws.Cells.Item(7, 1).Style = "Normal"

$ws.Cells.Item(7, 2).Value = "حسن "
$ws.Cells.Item(7, 2).Style = "Normal"

# C7 ("2222") looks numeric, so force text storage the same way the rest of
# the column already is, then restore Normal style so no formatting sticks.
$ws.Cells.Item(7, 3).NumberFormat = "@"
$ws.Cells.Item(7, 3).Value = "2222"
$ws.Cells.Item(7, 3).Style = "Normal"

$ws.Cells.Item(7, 4).Value = "ايتا"
$ws.Cells.Item(7, 4).Style = "Normal"

$ws.Cells.Item(7, 5).Value = "الرحلة 2"
$ws.Cells.Item(7, 5).Style = "Normal"

$ws.Cells.Item(7, 6).Value = "C3"
$ws.Cells.Item(7, 6).Style = "Normal"

$ws.Cells.Item(7, 7).Value = "NRC"
$ws.Cells.Item(7, 7).Style = "Normal"

$ws.Cells.Item(7, 8).Value = "٠٢‏/٠٥‏/٢٠٢٥ ٠٢:٠٠:٣٥ م"
$ws.Cells.Item(7, 8).Style = "Normal"
